$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Workbook-level: change active tab to "Regista Apartamento" (index 8, 0-based)
#    and move the tabSelected flag from "Autenticar Utilizador" to "Regista Apartamento"
# ---------------------------------------------------------------------------
$wsReg = $wb.Worksheets.Item("Regista Apartamento")
$wsAut = $wb.Worksheets.Item("Autenticar Utilizador")

# ---------------------------------------------------------------------------
# 2. Restructure the "Regista Apartamento" sheet
# ---------------------------------------------------------------------------

# Insert 5 new rows (9-13) for the new "Senhorio" related flow-of-events steps
$wsReg.Range("A9:A13").EntireRow.Insert()
# Insert 7 new rows (19-25) for the new exception block + spacing before the
# (now pushed down) "dados do apartamento invalidos" exception block
$wsReg.Range("A19:A25").EntireRow.Insert()

# --- Flow of events (rows 8-16) ---
$wsReg.Range("D8").Value = "Verifica se existe um Senhorio no sistema"

# --- Preconditions text change (row 5) ---
$wsReg.Range("C5").Value = "Não existe um senhorio nem um apartamento no sistema"

$wsReg.Range("C9").Value = "Fornece dados do Senhorio"
$wsReg.Range("B9").Formula = "=1+B8"

$wsReg.Range("D10").Value = "Valida dados do Senhorio"
$wsReg.Range("B10:B16").Formula = "=1+B9"

$wsReg.Range("D11").Value = "Regista Senhorio"
$wsReg.Range("D12").Value = "Indica que o Senhorio foi registado"

$wsReg.Range("C13").Value = "Fornece dados do apartamento"
$wsReg.Range("D14").Value = "Valida dados do apartamento"
$wsReg.Range("D15").Value = "Regista o apartamento"
$wsReg.Range("D16").Value = "Indica que o apartamento foi registado"

# --- New exception block 1 : "ja existe um Senhorio" (rows 17-18) ---
$wsReg.Range("C17").Value = "Actor Input"
$wsReg.Range("D17").Value = "System Response"
$wsReg.Range("A18").Value = "Excepção 1               (passo 1)" + [char]10 + "[existe um Senhorio no sistema]"
$wsReg.Range("B18").Value = 1
$wsReg.Range("D18").Value = "Indica que já existe um Senhorio no Sistema"

# --- New exception block 2 : "dados do Senhorio invalidos" (rows 19-21) ---
$wsReg.Range("C19").Value = "Actor Input"
$wsReg.Range("D19").Value = "System Response"
$wsReg.Range("A20").Value = "Excepção 2               (passo 3)" + [char]10 + "[dados do Senhorio inválidos]"
$wsReg.Range("B20").Value = 1
$wsReg.Range("D20").Value = "Indica que os dados inseridos são inválidos"
$wsReg.Range("B21").Value = 2
$wsReg.Range("D21").Value = "Regressa a 2"

# --- Existing exception block, now shifted to rows 26-28, text tweaks ---
$wsReg.Range("D28").Value = "Regressa a 6"
$wsReg.Range("A27").Value = "Excepção 3               (passo 7)" + [char]10 + "[dados do apartamento inválidos]"

# --- Formatting: bold + border + center + wrap for the two new "Actor Input / System Response" header rows ---
foreach ($r in @(17,19)) {
    $rng = $wsReg.Range("A" + $r + ":D" + $r)
    $rng.Font().Bold = $true
    $rng.Borders().LineStyle = 1
    $rng.HorizontalAlignment = -4108
    $rng.VerticalAlignment = -4108
}

# --- Formatting: bold + border + center + wrap text for the exception column A cells ---
foreach ($r in @(18,20,21)) {
    $rng = $wsReg.Range("A" + $r)
    $rng.Font().Bold = $true
    $rng.Borders().LineStyle = 1
    $rng.HorizontalAlignment = -4108
    $rng.VerticalAlignment = -4108
    $rng.WrapText = $true
}

# --- Formatting: bordered/centered/wrap cells for B/C/D in the exception rows ---
foreach ($r in @(18,20,21)) {
    $rng = $wsReg.Range("B" + $r + ":D" + $r)
    $rng.Borders().LineStyle = 1
    $rng.HorizontalAlignment = -4108
    $rng.VerticalAlignment = -4108
}
foreach ($r in @(18,20)) {
    $rng = $wsReg.Range("D" + $r)
    $rng.WrapText = $true
}

# Row heights for the wrapped exception rows
$wsReg.Rows.Item(18).RowHeight = 72
$wsReg.Rows.Item(20).RowHeight = 72

# --- Column widths ---
$wsReg.Columns.Item(1).ColumnWidth = 18.71
$wsReg.Columns.Item(3).ColumnWidth = 46.29

# --- Sheet view: this sheet becomes the selected/active one ---
$wsAut.Activate()
$wsReg.Activate()
$wsReg.Application.ActiveWindow.ScrollRow = 13
$wsReg.Range("A28").Select()

$wb.Worksheets.Item("Regista Apartamento").Activate()
